$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows 2,3,4,5,6,7,9 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 451
$wsExhibit.Range("F3").Value = 34
$wsExhibit.Range("F4").Value = 66
$wsExhibit.Range("F5").Value = 5104
$wsExhibit.Range("F6").Value = 173
$wsExhibit.Range("F7").Value = 36
$wsExhibit.Range("F9").Value = 317

# Sheet "全部类型" (all types) - rows 2,7,8,9,10,11,14 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 451
$wsAll.Range("F7").Value = 34
$wsAll.Range("F8").Value = 66
$wsAll.Range("F9").Value = 5104
$wsAll.Range("F10").Value = 173
$wsAll.Range("F11").Value = 36
$wsAll.Range("F14").Value = 317
